{"js": "// Update the date line and the twenty-six \"NN\u00d7NN=\" multiplication\n// prompts in the practice-sheet table. Every old value is unique in\n// the document, so a direct search + full-text replace per pair is\n// safe and unambiguous.\nconst replacements = [\n  [\"2024-03-28 Thursday\", \"2024-03-29 Friday\"],\n  [\"52\u00d720=\", \"89\u00d788=\"],\n  [\"51\u00d778=\", \"60\u00d785=\"],\n  [\"47\u00d785=\", \"95\u00d757=\"],\n  [\"27\u00d711=\", \"57\u00d741=\"],\n  [\"18\u00d796=\", \"41\u00d753=\"],\n  [\"51\u00d727=\", \"69\u00d737=\"],\n  [\"70\u00d776=\", \"17\u00d774=\"],\n  [\"41\u00d742=\", \"52\u00d744=\"],\n  [\"61\u00d741=\", \"63\u00d741=\"],\n  [\"53\u00d778=\", \"67\u00d766=\"],\n  [\"90\u00d778=\", \"16\u00d729=\"],\n  [\"30\u00d733=\", \"96\u00d731=\"],\n  [\"61\u00d777=\", \"59\u00d728=\"],\n  [\"22\u00d741=\", \"44\u00d781=\"],\n  [\"34\u00d747=\", \"24\u00d785=\"],\n  [\"49\u00d787=\", \"55\u00d711=\"],\n  [\"17\u00d766=\", \"72\u00d774=\"],\n  [\"58\u00d717=\", \"54\u00d744=\"],\n  [\"71\u00d715=\", \"69\u00d752=\"],\n  [\"39\u00d732=\", \"71\u00d734=\"],\n  [\"29\u00d716=\", \"14\u00d774=\"],\n  [\"98\u00d775=\", \"50\u00d799=\"],\n  [\"83\u00d720=\", \"93\u00d711=\"],\n  [\"77\u00d719=\", \"78\u00d772=\"],\n  [\"19\u00d747=\", \"66\u00d738=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-six \"NN\u00d7NN=\" multiplication\n# prompts in the practice-sheet table. Every old value is unique in\n# the document, so Find/Replace (ReplaceAll) per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-03-28 Thursday\"; New = \"2024-03-29 Friday\" },\n    @{ Old = \"52\u00d720=\";              New = \"89\u00d788=\" },\n    @{ Old = \"51\u00d778=\";              New = \"60\u00d785=\" },\n    @{ Old = \"47\u00d785=\";              New = \"95\u00d757=\" },\n    @{ Old = \"27\u00d711=\";              New = \"57\u00d741=\" },\n    @{ Old = \"18\u00d796=\";              New = \"41\u00d753=\" },\n    @{ Old = \"51\u00d727=\";              New = \"69\u00d737=\" },\n    @{ Old = \"70\u00d776=\";              New = \"17\u00d774=\" },\n    @{ Old = \"41\u00d742=\";              New = \"52\u00d744=\" },\n    @{ Old = \"61\u00d741=\";              New = \"63\u00d741=\" },\n    @{ Old = \"53\u00d778=\";              New = \"67\u00d766=\" },\n    @{ Old = \"90\u00d778=\";              New = \"16\u00d729=\" },\n    @{ Old = \"30\u00d733=\";              New = \"96\u00d731=\" },\n    @{ Old = \"61\u00d777=\";              New = \"59\u00d728=\" },\n    @{ Old = \"22\u00d741=\";              New = \"44\u00d781=\" },\n    @{ Old = \"34\u00d747=\";              New = \"24\u00d785=\" },\n    @{ Old = \"49\u00d787=\";              New = \"55\u00d711=\" },\n    @{ Old = \"17\u00d766=\";              New = \"72\u00d774=\" },\n    @{ Old = \"58\u00d717=\";              New = \"54\u00d744=\" },\n    @{ Old = \"71\u00d715=\";              New = \"69\u00d752=\" },\n    @{ Old = \"39\u00d732=\";              New = \"71\u00d734=\" },\n    @{ Old = \"29\u00d716=\";              New = \"14\u00d774=\" },\n    @{ Old = \"98\u00d775=\";              New = \"50\u00d799=\" },\n    @{ Old = \"83\u00d720=\";              New = \"93\u00d711=\" },\n    @{ Old = \"77\u00d719=\";              New = \"78\u00d772=\" },\n    @{ Old = \"19\u00d747=\";              New = \"66\u00d738=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
